$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A2").Value2 = "ECs"
$ws.Range("B2").Value2 = "C3"
$ws.Range("C2").Value2 = "Itgam"
$ws.Range("D2").Value2 = "ECs"
$ws.Range("E2").Value2 = 3
$ws.Range("F2").Value2 = 1
$ws.Range("G2").Value2 = 40.942832
$ws.Range("H2").Value2 = 122.828496
$ws.Range("I2").Value2 = 0.2583000005785167
$ws.Range("J2").Value2 = 0.2583000005785167
$ws.Range("K2").Value2 = 1
$ws.Range("L2").Value2 = 0.3333333333333333
$ws.Range("M2").Value2 = 66.38494866666666
$ws.Range("N2").Value2 = 199.154846
$ws.Range("O2").Value2 = 0.4171200956172241
$ws.Range("P2").Value2 = 0.4171200956172242
$ws.Range("Q2").Value2 = 2717.987800587957
$ws.Range("R2").Value2 = 24461.89020529162
$ws.Range("S2").Value2 = 0.1077421209392399
$ws.Range("T2").Value2 = 0.1077421209392399
$ws.Range("A3").Value2 = "ECs"
$ws.Range("B3").Value2 = "C3"
$ws.Range("C3").Value2 = "Itgam"
$ws.Range("D3").Value2 = "M2"
$ws.Range("E3").Value2 = 3
$ws.Range("F3").Value2 = 1
$ws.Range("G3").Value2 = 40.942832
$ws.Range("H3").Value2 = 122.828496
$ws.Range("I3").Value2 = 0.2583000005785167
$ws.Range("J3").Value2 = 0.2583000005785167
$ws.Range("K3").Value2 = 3
$ws.Range("L3").Value2 = 1
$ws.Range("M3").Value2 = 92.76573566666667
$ws.Range("N3").Value2 = 278.297207
$ws.Range("O3").Value2 = 0.5828799043827758
$ws.Range("P3").Value2 = 0.5828799043827758
$ws.Range("Q3").Value2 = 3798.091930756742
$ws.Range("R3").Value2 = 34182.82737681067
$ws.Range("S3").Value2 = 0.1505578796392767
$ws.Range("T3").Value2 = 0.1505578796392767
$ws.Range("A4").Value2 = "FAPs"
$ws.Range("B4").Value2 = "C3"
$ws.Range("C4").Value2 = "Itgam"
$ws.Range("D4").Value2 = "ECs"
$ws.Range("E4").Value2 = 3
$ws.Range("F4").Value2 = 1
$ws.Range("G4").Value2 = 79.68771233333334
$ws.Range("H4").Value2 = 239.063137
$ws.Range("I4").Value2 = 0.5027335710876245
$ws.Range("J4").Value2 = 0.5027335710876245
$ws.Range("K4").Value2 = 1
$ws.Range("L4").Value2 = 0.3333333333333333
$ws.Range("M4").Value2 = 66.38494866666666
$ws.Range("N4").Value2 = 199.154846
$ws.Range("O4").Value2 = 0.4171200956172241
$ws.Range("P4").Value2 = 0.4171200956172242
$ws.Range("Q4").Value2 = 5290.064692612434
$ws.Range("R4").Value2 = 47610.5822335119
$ws.Range("S4").Value2 = 0.2097002752420585
$ws.Range("T4").Value2 = 0.2097002752420585
$ws.Range("A5").Value2 = "FAPs"
$ws.Range("B5").Value2 = "C3"
$ws.Range("C5").Value2 = "Itgam"
$ws.Range("D5").Value2 = "M2"
$ws.Range("E5").Value2 = 3
$ws.Range("F5").Value2 = 1
$ws.Range("G5").Value2 = 79.68771233333334
$ws.Range("H5").Value2 = 239.063137
$ws.Range("I5").Value2 = 0.5027335710876245
$ws.Range("J5").Value2 = 0.5027335710876245
$ws.Range("K5").Value2 = 3
$ws.Range("L5").Value2 = 1
$ws.Range("M5").Value2 = 92.76573566666667
$ws.Range("N5").Value2 = 278.297207
$ws.Range("O5").Value2 = 0.5828799043827758
$ws.Range("P5").Value2 = 0.5828799043827758
$ws.Range("Q5").Value2 = 7392.289258195374
$ws.Range("R5").Value2 = 66530.60332375836
$ws.Range("S5").Value2 = 0.293033295845566
$ws.Range("T5").Value2 = 0.293033295845566
$ws.Range("A6").Value2 = "M2"
$ws.Range("B6").Value2 = "C3"
$ws.Range("C6").Value2 = "Itgam"
$ws.Range("D6").Value2 = "ECs"
$ws.Range("E6").Value2 = 3
$ws.Range("F6").Value2 = 1
$ws.Range("G6").Value2 = 37.53186833333334
$ws.Range("H6").Value2 = 112.595605
$ws.Range("I6").Value2 = 0.2367809244903433
$ws.Range("J6").Value2 = 0.2367809244903433
$ws.Range("K6").Value2 = 1
$ws.Range("L6").Value2 = 0.3333333333333333
$ws.Range("M6").Value2 = 66.38494866666666
$ws.Range("N6").Value2 = 199.154846
$ws.Range("O6").Value2 = 0.4171200956172241
$ws.Range("P6").Value2 = 0.4171200956172242
$ws.Range("Q6").Value2 = 2491.551152672425
$ws.Range("R6").Value2 = 22423.96037405183
$ws.Range("S6").Value2 = 0.09876608186374672
$ws.Range("T6").Value2 = 0.09876608186374672
$ws.Range("A7").Value2 = "M2"
$ws.Range("B7").Value2 = "C3"
$ws.Range("C7").Value2 = "Itgam"
$ws.Range("D7").Value2 = "M2"
$ws.Range("E7").Value2 = 3
$ws.Range("F7").Value2 = 1
$ws.Range("G7").Value2 = 37.53186833333334
$ws.Range("H7").Value2 = 112.595605
$ws.Range("I7").Value2 = 0.2367809244903433
$ws.Range("J7").Value2 = 0.2367809244903433
$ws.Range("K7").Value2 = 3
$ws.Range("L7").Value2 = 1
$ws.Range("M7").Value2 = 92.76573566666667
$ws.Range("N7").Value2 = 278.297207
$ws.Range("O7").Value2 = 0.5828799043827758
$ws.Range("P7").Value2 = 0.5828799043827758
$ws.Range("Q7").Value2 = 3481.671376886138
$ws.Range("R7").Value2 = 31335.04239197524
$ws.Range("S7").Value2 = 0.1380148426265966
$ws.Range("T7").Value2 = 0.1380148426265966
$ws.Range("A8").Value2 = "sCs"
$ws.Range("B8").Value2 = "C3"
$ws.Range("C8").Value2 = "Itgam"
$ws.Range("D8").Value2 = "ECs"
$ws.Range("E8").Value2 = 3
$ws.Range("F8").Value2 = 1
$ws.Range("G8").Value2 = 0.3464216666666666
$ws.Range("H8").Value2 = 1.039265
$ws.Range("I8").Value2 = 0.002185503843515531
$ws.Range("J8").Value2 = 0.002185503843515531
$ws.Range("K8").Value2 = 1
$ws.Range("L8").Value2 = 0.3333333333333333
$ws.Range("M8").Value2 = 66.38494866666666
$ws.Range("N8").Value2 = 199.154846
$ws.Range("O8").Value2 = 0.4171200956172241
$ws.Range("P8").Value2 = 0.4171200956172242
$ws.Range("Q8").Value2 = 22.99718455868777
$ws.Range("R8").Value2 = 206.97466102819
$ws.Range("S8").Value2 = 0.0009116175721790093
$ws.Range("T8").Value2 = 0.0009116175721790095
$ws.Range("A9").Value2 = "sCs"
$ws.Range("B9").Value2 = "C3"
$ws.Range("C9").Value2 = "Itgam"
$ws.Range("D9").Value2 = "M2"
$ws.Range("E9").Value2 = 3
$ws.Range("F9").Value2 = 1
$ws.Range("G9").Value2 = 0.3464216666666666
$ws.Range("H9").Value2 = 1.039265
$ws.Range("I9").Value2 = 0.002185503843515531
$ws.Range("J9").Value2 = 0.002185503843515531
$ws.Range("K9").Value2 = 3
$ws.Range("L9").Value2 = 1
$ws.Range("M9").Value2 = 92.76573566666667
$ws.Range("N9").Value2 = 278.297207
$ws.Range("O9").Value2 = 0.5828799043827758
$ws.Range("P9").Value2 = 0.5828799043827758
$ws.Range("Q9").Value2 = 32.13606075920611
$ws.Range("R9").Value2 = 289.224546832855
$ws.Range("S9").Value2 = 0.001273886271336522
$ws.Range("T9").Value2 = 0.001273886271336522
